$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 635
$ws1.Range("F3").Value = 2213
$ws1.Range("F5").Value = 13235
$ws1.Range("F12").Value = 13792
$ws1.Range("F13").Value = 14418
$ws1.Range("F23").Value = 113
$ws1.Range("F25").Value = 5470
$ws1.Range("F27").Value = 376

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 635
$ws4.Range("F3").Value = 2213
$ws4.Range("F5").Value = 13235
$ws4.Range("F13").Value = 13792
$ws4.Range("F14").Value = 14418
$ws4.Range("F24").Value = 113
$ws4.Range("F26").Value = 5470
$ws4.Range("F27").Value = 941
$ws4.Range("F28").Value = 376
